$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was an empty placeholder row; clear it so it drops out of the saved
# sheetData entirely (the "Marcelo"/101 record stays at row 3, unmoved).
$ws.Range("A2:C2").ClearContents()

# New habitacion/paciente record, added as row 4.
$ws.Range("A4").Value = "Juan Perez"
# Leading apostrophe forces "205" to be stored as text (matching how the
# adjacent room-number cell B3 is stored) instead of being auto-coerced to
# a number.
$ws.Range("B4").Value = "'205"
